# Update column F ("dSF") values on Sheet1 to reflect the repulled data /
# recalculated means, per the commit "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = 3
    9  = 0
    12 = -1
    15 = -2
    16 = 0
    20 = -1
    31 = 0
    32 = 0
    34 = 0
    35 = 2
    36 = 0
    48 = 2
    50 = 1
    58 = 3
    60 = -2
    64 = -7
    67 = 5
    73 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}

$wb.Save()
